$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 25.91000000000061
$ws.Range("H2").Value = [double]"1.15347846714302e-15"
$ws.Range("K2").Value = 35.94130627491537
$ws.Range("L2").Value = "[28.339492578758453, 43.54311997107229]"
$ws.Range("O2").Value = 1.603816069400195
$ws.Range("P2").Value = "[1.3648160276856558, 1.8428161111147343]"
$ws.Range("S2").Value = 58.11143102589541
$ws.Range("T2").Value = "[52.980438232097086, 63.24242381969374]"
$ws.Range("W2").Value = 19.29633633633679
$ws.Range("X2").Value = 18.3107707707712
$ws.Range("Y2").Value = 20.28190190190238

# Row 3
$ws.Range("E3").Value = 23.70000000000027
$ws.Range("G3").Value = [double]"2.594127271215996e-06"
$ws.Range("H3").Value = [double]"8.782235579943773e-06"
$ws.Range("K3").Value = 32.1059606005995
$ws.Range("L3").Value = "[15.826883413478793, 48.38503778772021]"
$ws.Range("M3").Value = 0.0001488923094397077
$ws.Range("N3").Value = 0.0001488923094397077
$ws.Range("O3").Value = -2.226474072814388
$ws.Range("P3").Value = "[-2.7170531058073886, -1.7358950398213873]"
$ws.Range("Q3").Value = [double]"1.554312234475219e-15"
$ws.Range("R3").Value = [double]"1.554312234475219e-15"
$ws.Range("S3").Value = 65.24693200733429
$ws.Range("T3").Value = "[56.856767138543134, 73.63709687612545]"
$ws.Range("W3").Value = 8.398198198198292
$ws.Range("X3").Value = 6.547747747747821
$ws.Range("Y3").Value = 10.24864864864876
